$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that was bumped
# from 45186 (2023-09-17) to 45188 (2023-09-19) for every data row
# (rows 2 through 472).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

$range = $ws.Range("C2:C$lastRow")
$range.Value = 45188
